$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the current data rows (2-12), columns A-D, into memory.
$rows = @()
for ($r = 2; $r -le 12; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $d = $ws.Cells.Item($r, 4).Value()
    $rows += ,@($a, $b, $c, $d)
}

# Sort the rows ascending by column A (time).
$sorted = $rows | Sort-Object { $_[0] }

# Write the sorted rows back to the sheet.
for ($i = 0; $i -lt $sorted.Count; $i++) {
    $r = $i + 2
    $row = $sorted[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
